$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.987.30'
$ws.Range('E2').Value = '  +1.99%  '
$ws.Range('D3').Value = '2.306.78'
$ws.Range('E3').Value = '  +1.72%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '''302.34'
$ws.Range('E5').Value = '  +1.07%  '
$ws.Range('D6').Value = '''100.57'
$ws.Range('E6').Value = '  +5.08%  '
$ws.Range('E7').Value = '  +1.65%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '''0.518'
$ws.Range('E9').Value = '  +5.17%  '
$ws.Range('D10').Value = '''35.05'
$ws.Range('E10').Value = '  +5.32%  '
$ws.Range('D11').Value = '''0.0795'
$ws.Range('E11').Value = '  +0.91%  '
$ws.Range('E12').Value = '  +3.97%  '
$ws.Range('D13').Value = '''17.86'
$ws.Range('E13').Value = '  +15.24%  '
$ws.Range('D14').Value = '''6.90'
$ws.Range('E14').Value = '  +3.48%  '
$ws.Range('D15').Value = '2.683.33'
$ws.Range('E15').Value = '  +2.50%  '
$ws.Range('D16').Value = '2.313.95'
$ws.Range('E16').Value = '  +2.22%  '
$ws.Range('D17').Value = '''0.811'
$ws.Range('E17').Value = '  +3.63%  '
$ws.Range('D18').Value = '42.900.91'
$ws.Range('E18').Value = '  +1.98%  '
$ws.Range('D19').Value = '''12.61'
$ws.Range('E19').Value = '  +8.57%  '
$ws.Range('E20').Value = '  +3.25%  '
$ws.Range('D22').Value = '''67.92'
$ws.Range('E22').Value = '  +2.24%  '
$ws.Range('D23').Value = '''237.03'
$ws.Range('E23').Value = '  +0.87%  '
$ws.Range('D24').Value = '''2.20'
$ws.Range('E24').Value = '  +11.95%  '
$ws.Range('D25').Value = '''2.46'
$ws.Range('E25').Value = '  +0.59%  '
$ws.Range('D26').Value = '''0.998'
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('D27').Value = '''24.72'
$ws.Range('E27').Value = '  +3.09%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').Value = '''167.56'
$ws.Range('E28').Value = '  -0.45%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '''2.08'
$ws.Range('E29').Value = '  -4.46%  '
$ws.Range('D30').Value = '''33.92'
$ws.Range('E30').Value = '  +1.13%  '
$ws.Range('D31').Value = '''9.21'
$ws.Range('E31').Value = '  +0.62%  '
$ws.Range('E32').Value = '  +0.12%  '
$ws.Range('D33').Value = '''5.01'
$ws.Range('E33').Value = '  +2.43%  '
$ws.Range('E34').Value = '  +2.89%  '
$ws.Range('E35').Value = '  +3.65%  '
$ws.Range('D36').Value = '''16.93'
$ws.Range('E36').Value = '  +2.47%  '
$ws.Range('D37').Value = '''0.0690'
$ws.Range('E37').Value = '  +0.63%  '
$ws.Range('E38').Value = '  +3.46%  '
$ws.Range('E39').Value = '  +1.75%  '
$ws.Range('E40').Value = '  +3.69%  '
$ws.Range('E41').Value = '  +0.86%  '
$ws.Range('B42').Value = 'ApeXProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D42').Value = '''2.30'
$ws.Range('E42').Value = '  -5.19%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.000.74'
$ws.Range('E43').Value = '  +1.76%  '
$ws.Range('D44').Value = '''0.0287'
$ws.Range('E44').Value = '  +3.87%  '
$ws.Range('D45').Value = '''10.24'
$ws.Range('E45').Value = '  +7.56%  '
$ws.Range('D46').Value = '''17.53'
$ws.Range('E46').Value = '  +1.46%  '
$ws.Range('E47').Value = '  +2.59%  '
$ws.Range('D48').Value = '''55.72'
$ws.Range('E48').Value = '  +6.90%  '
$ws.Range('D49').Value = '2.526.99'
$ws.Range('E49').Value = '  +1.43%  '
$ws.Range('D50').Value = '''1.52'
$ws.Range('E50').Value = '  +4.17%  '
$ws.Range('D51').Value = '''4.55'
$ws.Range('E51').Value = '  +1.18%  '
